$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was added to the "Jengibre" (ginger) series.
# It belongs chronologically between the current row 15 (2021-03-29) and the
# current row 16 (2021-04-29), so insert a fresh row at position 16, which
# pushes the former rows 16..51 down to 17..52.
$ws.Rows.Item(16).Insert()

# Populate the newly inserted row 16 with the new observation's data.
$ws.Cells.Item(16, 1).Value = 6
$ws.Cells.Item(16, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(16, 3).Value = "Metropolitana"
$ws.Cells.Item(16, 4).Value = 44494
$ws.Cells.Item(16, 5).Value = 13
$ws.Cells.Item(16, 6).Value = 100114007
$ws.Cells.Item(16, 7).Value = "Jengibre"
$ws.Cells.Item(16, 8).Value = "Sin especificar"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 370
$ws.Cells.Item(16, 11).Value = 13000
$ws.Cells.Item(16, 12).Value = 15000
$ws.Cells.Item(16, 13).Value = 14243
$ws.Cells.Item(16, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(16, 15).Value = "Perú"
$ws.Cells.Item(16, 16).Value = 1096
$ws.Cells.Item(16, 17).Value = 13
$ws.Cells.Item(16, 18).Value = "Hortaliza"
